# HOTFIX change date format in example
#
# The example import sheet had six "Дата *" columns (K:P, "Дата поверки",
# "Дата следующей поверки", "Дата установки", "Дата ввода в эксплуатацию",
# "Дата опломбирования", "Дата контрольных показаний") whose sample values
# used a dotted dd.mm.yyyy format (and varied row to row). Excel's importer
# expects ISO (yyyy-mm-dd) dates, so every data row's K:P sample values are
# rewritten to the same canonical set of ISO date strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("2021-01-20", "2021-01-21", "2021-01-22", "2021-01-23", "2021-01-24", "2021-01-25")

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

# Header is row 1; data starts on row 2.
for ($row = 2; $row -le $lastRow; $row++) {
    for ($i = 0; $i -lt $dates.Length; $i++) {
        $col = 11 + $i  # K=11, L=12, M=13, N=14, O=15, P=16
        $ws.Cells.Item($row, $col).Value = $dates[$i]
    }
}
